# Apply cryptos list update (price + volume refresh) from GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a numeric-looking string while preserving it as literal text
# (matches source data which stores prices as inline text, not numbers).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2: Bitcoin
$ws.Range("D2").Value = '39.436.61'
$ws.Range("E2").Value = '  +1.78%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.157.79'
$ws.Range("E3").Value = '  +2.82%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.08%  '

# Row 5: BNB
Set-TextValue $ws.Range("D5") '227.68'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6: XRP
$ws.Range("E6").Value = '  +0.89%  '

# Row 7: Solana
Set-TextValue $ws.Range("D7") '63.98'
$ws.Range("E7").Value = '  +3.98%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.00%  '

# Row 9: Cardano
Set-TextValue $ws.Range("D9") '0.396'
$ws.Range("E9").Value = '  +2.50%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +1.37%  '

# Row 11: TRON
$ws.Range("E11").Value = '  +0.49%  '

# Row 12: Chainlink
Set-TextValue $ws.Range("D12") '15.99'
$ws.Range("E12").Value = '  +2.99%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.478.17'
$ws.Range("E13").Value = '  +2.72%  '

# Row 14: Avalanche
Set-TextValue $ws.Range("D14") '22.11'
$ws.Range("E14").Value = '  +0.42%  '

# Row 15: Polygon
Set-TextValue $ws.Range("D15") '0.812'
$ws.Range("E15").Value = '  +0.55%  '

# Row 16: Polkadot
$ws.Range("E16").Value = '  +0.80%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.158.17'
$ws.Range("E17").Value = '  +3.11%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '39.372.92'
$ws.Range("E18").Value = '  +1.50%  '

# Row 19: Litecoin
Set-TextValue $ws.Range("D19") '71.80'
$ws.Range("E19").Value = '  -0.27%  '

# Row 20: Uniswap
Set-TextValue $ws.Range("D20") '6.10'

# Row 21: ShibaInu
$ws.Range("E21").Value = '  +1.26%  '

# Row 22: BitcoinCash
Set-TextValue $ws.Range("D22") '230.81'
$ws.Range("E22").Value = '  +1.41%  '

# Row 23: Dai
$ws.Range("E23").Value = '  +0.00%  '

# Row 24: PancakeSwap -> Toncoin
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D24") '2.34'
$ws.Range("E24").Value = '  -2.04%  '

# Row 25: Toncoin -> PancakeSwap
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D25") '2.35'
$ws.Range("E25").Value = '  +0.31%  '

# Row 26: Monero
Set-TextValue $ws.Range("D26") '172.13'
$ws.Range("E26").Value = '  +0.22%  '

# Row 27: Cosmos
Set-TextValue $ws.Range("D27") '9.49'
$ws.Range("E27").Value = '  -0.56%  '

# Row 28: Kaspa
Set-TextValue $ws.Range("D28") '0.140'
$ws.Range("E28").Value = '  +1.69%  '

# Row 29: EthereumClassic
$ws.Range("E29").Value = '  +2.83%  '

# Row 30: ImmutableX
Set-TextValue $ws.Range("D30") '1.41'
$ws.Range("E30").Value = '  -0.32%  '

# Row 31: WEMIXToken
$ws.Range("E31").Value = '  +7.89%  '

# Row 32: Stellar
$ws.Range("E32").Value = '  +0.61%  '

# Row 33: Filecoin
Set-TextValue $ws.Range("D33") '4.60'
$ws.Range("E33").Value = '  +1.51%  '

# Row 34: THORChain
$ws.Range("E34").Value = '  +9.24%  '

# Row 35: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D35") '4.73'
$ws.Range("E35").Value = '  -0.57%  '

# Row 36: Hedera
$ws.Range("E36").Value = '  -0.75%  '

# Row 37: LidoDAOToken
$ws.Range("E37").Value = '  +0.34%  '

# Row 38: RenderToken
Set-TextValue $ws.Range("D38") '3.57'
$ws.Range("E38").Value = '  +0.19%  '

# Row 40: Aave
Set-TextValue $ws.Range("D40") '103.60'
$ws.Range("E40").Value = '  +2.19%  '

# Row 41: VeChain
$ws.Range("E41").Value = '  +0.60%  '

# Row 42: InjectiveProtocol -> Maker
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.540.25'
$ws.Range("E42").Value = '  +0.38%  '

# Row 43: Maker -> InjectiveProtocol
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D43") '17.60'
$ws.Range("E43").Value = '  -3.21%  '

# Row 44: TrustWalletToken
$ws.Range("E44").Value = '  +4.31%  '

# Row 45: FTXToken -> Cronos
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D45") '0.0931'
$ws.Range("E45").Value = '  +2.35%  '

# Row 46: Cronos -> FTXToken
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D46") '4.31'
$ws.Range("E46").Value = '  +4.58%  '

# Row 47: HuobiToken
$ws.Range("E47").Value = '  +0.65%  '

# Row 48: ARBITRUM
Set-TextValue $ws.Range("D48") '1.09'
$ws.Range("E48").Value = '  +4.93%  '

# Row 49: FraxShare
Set-TextValue $ws.Range("D49") '7.69'
$ws.Range("E49").Value = '  -0.89%  '

# Row 50: RocketPoolETH
$ws.Range("D50").Value = '2.361.62'
$ws.Range("E50").Value = '  +2.88%  '

# Row 51: MXToken
$ws.Range("E51").Value = '  +0.05%  '
